# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates columns I (DAMSLTag) and J (DialogAct) for specific rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    8   = @{ I = "ba"; J = "Appreciation" }
    9   = @{ I = "sd"; J = "Statement-non-opinion" }
    12  = @{ I = "sd"; J = "Statement-non-opinion" }
    14  = @{ I = "sv"; J = "Statement-opinion" }
    15  = @{ I = "sv"; J = "Statement-opinion" }
    17  = @{ I = "sv"; J = "Statement-opinion" }
    18  = @{ I = "sd"; J = "Statement-non-opinion" }
    19  = @{ I = "sd"; J = "Statement-non-opinion" }
    20  = @{ I = "sv"; J = "Statement-opinion" }
    43  = @{ I = "b";  J = "Acknowledge (Backchannel)" }
    60  = @{ I = "ba"; J = "Appreciation" }
    78  = @{ I = "b";  J = "Acknowledge (Backchannel)" }
    81  = @{ I = "sd"; J = "Statement-non-opinion" }
    83  = @{ I = "sd"; J = "Statement-non-opinion" }
    84  = @{ I = "sd"; J = "Statement-non-opinion" }
    93  = @{ I = "sv"; J = "Statement-opinion" }
    95  = @{ I = "%";  J = "Uninterpretable" }
    99  = @{ I = "sd"; J = "Statement-non-opinion" }
    103 = @{ I = "ba"; J = "Appreciation" }
    111 = @{ I = "sd"; J = "Statement-non-opinion" }
    117 = @{ I = "sd"; J = "Statement-non-opinion" }
}

foreach ($rowNum in $updates.Keys) {
    $vals = $updates[$rowNum]
    $ws.Range("I$rowNum").Value = $vals.I
    $ws.Range("J$rowNum").Value = $vals.J
}
